$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("textbook")
$rng.Collapse(0)
$rng.InsertAfter("s")
$rng2 = $d.Range($rng.Start, $rng.Start+1)
$rng2.Bold = 1
$rng2.Bold = 0
# Now test: does a subsequent no-op find/replace merge them back?
$d.Content.Find.Execute("two textbooks and crawled", $true, $false, $false, $false, $false, $true, 1, $false, "two textbooks and crawled", 2)
